# Rename the "series A/B/C/D" / "Other series" / "Over 65" / "Under 66"
# labels to the new "type A/B/C/D" / "Other types" / "> 65" / "< 66" labels
# on both worksheets, and update the remembered cell selection on each
# sheet's view (as recorded in the xlsx after the edit).

$wb = $excel.ActiveWorkbook

# --- Sheet "metadata_included" ---------------------------------------
$ws1 = $wb.Worksheets.Item("metadata_included")

$ws1.Range("D3").Value = "...Other. types ….........."
$ws1.Range("E3").Value = "type               A"
$ws1.Range("F3").Value = "type B"
$ws1.Range("G3").Value = "type C"
$ws1.Range("H3").Value = "type D"

$ws1.Range("C4").Value = "< 66"
$ws1.Range("C5").Value = "> 65"

$ws1.Activate() | Out-Null
$ws1.Range("D3:H3").Select() | Out-Null

# --- Sheet "no_metadata" ----------------------------------------------
$ws2 = $wb.Worksheets.Item("no_metadata")

$ws2.Range("D1").Value = "...Other. types ….........."
$ws2.Range("E1").Value = "type               A"
$ws2.Range("F1").Value = "type B"
$ws2.Range("G1").Value = "type C"
$ws2.Range("H1").Value = "type D"

$ws2.Range("C2").Value = "< 66"
$ws2.Range("C3").Value = "> 65"

$ws2.Activate() | Out-Null
$ws2.Range("H9").Select() | Out-Null

